# Generate Report for Handoff
#
# Inserts a new tracked file "0cd4e0fa-534a-4c0b-bf1d-b9b57b1db9e7.md" (status
# "Ready for handoff") into the localization-status workbook. The new file is
# inserted as a new row right before the existing "984aa78d-...-md" /
# ".localization-config" rows (so it becomes row 6 on every sheet, pushing the
# two trailing rows down by one) on all three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$newBaseName   = "0cd4e0fa-534a-4c0b-bf1d-b9b57b1db9e7"
$newMdFile     = "$newBaseName.md"
$newHash       = "dcd3d332b0d7f0538e7c30bc8452a132e6c349a4"
$newZhHandoff  = "$newBaseName.$newHash.zh-cn.xlf"
$newDeHandoff  = "$newBaseName.$newHash.de-de.xlf"
$zhHandoffDate = "2016-02-26 04:55:01"
$deHandoffDate = "2016-02-26 04:55:11"
$readyStatus   = "Ready for handoff"
$epoch         = "0001-01-01 00:00:00"
$includeStatus = "Include"

$mdCommitHash  = "3f856d1a2d6a2b5a13a6a39cd2d3e1a0cfa7d8e1"
$zhCommitHash  = "a5c8b16d8c0fa3e1f4f9b2d6e7c5a3f1a9d2c4b0"
$deCommitHash  = "b7e3f29a5c1d8e6b0a4f2c9d7e5b3a1f6c8d4e2a"

function New-MdUrl($name) {
    return "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommitHash/e2e/$name"
}
function New-ZhHandoffUrl($name) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommitHash/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/$name"
}
function New-DeHandoffUrl($name) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommitHash/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/$name"
}

# ---------------------------------------------------------------------------
# Sheet "Overview": columns A (File Name), B (zh-cn status), C (de-de status)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Rows.Item(6).Insert()
$wsOverview.Range("A6").Value2 = $newMdFile
$wsOverview.Range("B6").Value2 = $readyStatus
$wsOverview.Range("C6").Value2 = $readyStatus

# Rebuild hyperlinks for column A (File Name) in row order.
$wsOverview.Cells.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), (New-MdUrl "20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.md"), "", "", "20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), (New-MdUrl "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md"), "", "", "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), (New-MdUrl "831315f2-da6a-4fc4-b117-60389f3c6074.md"), "", "", "831315f2-da6a-4fc4-b117-60389f3c6074.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), (New-MdUrl "ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.md"), "", "", "ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), (New-MdUrl $newMdFile), "", "", $newMdFile)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), (New-MdUrl "984aa78d-5405-411a-8b9f-4541e5a8b93a.md"), "", "", "984aa78d-5405-411a-8b9f-4541e5a8b93a.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A8"), (New-MdUrl ".localization-config"), "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "zh-cn": columns A (md file), B (status), C (handoff file), D (handoff
# datetime), G (handback datetime / placeholder), H (handoff reason)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows.Item(6).Insert()
$wsZh.Range("A6").Value2 = $newMdFile
$wsZh.Range("B6").Value2 = $readyStatus
$wsZh.Range("C6").Value2 = $newZhHandoff
$wsZh.Range("D6").Value2 = $zhHandoffDate
$wsZh.Range("G6").Value2 = $epoch
$wsZh.Range("H6").Value2 = $includeStatus

$wsZh.Cells.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), (New-MdUrl "20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.md"), "", "", "20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), (New-ZhHandoffUrl "20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.550666d79c2cc553a731e154650f5b8682684363.zh-cn.xlf"), "", "", "20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.550666d79c2cc553a731e154650f5b8682684363.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), (New-MdUrl "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md"), "", "", "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), (New-ZhHandoffUrl "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.zh-cn.xlf"), "", "", "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), (New-MdUrl "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md"), "", "", "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), (New-ZhHandoffUrl "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.zh-cn.xlf"), "", "", "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), (New-MdUrl "831315f2-da6a-4fc4-b117-60389f3c6074.md"), "", "", "831315f2-da6a-4fc4-b117-60389f3c6074.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), (New-ZhHandoffUrl "831315f2-da6a-4fc4-b117-60389f3c6074.33096a8b89da196732c33acf66abed619fc16792.zh-cn.xlf"), "", "", "831315f2-da6a-4fc4-b117-60389f3c6074.33096a8b89da196732c33acf66abed619fc16792.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), (New-MdUrl "ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.md"), "", "", "ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C5"), (New-ZhHandoffUrl "ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.fc1d393ccc271987bd0fbbef8d59cb1dfbb601f9.zh-cn.xlf"), "", "", "ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.fc1d393ccc271987bd0fbbef8d59cb1dfbb601f9.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), (New-MdUrl $newMdFile), "", "", $newMdFile)
$wsZh.Hyperlinks.Add($wsZh.Range("C6"), (New-ZhHandoffUrl $newZhHandoff), "", "", $newZhHandoff)
$wsZh.Hyperlinks.Add($wsZh.Range("A7"), (New-MdUrl "984aa78d-5405-411a-8b9f-4541e5a8b93a.md"), "", "", "984aa78d-5405-411a-8b9f-4541e5a8b93a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C7"), (New-ZhHandoffUrl "984aa78d-5405-411a-8b9f-4541e5a8b93a.38149897969e365477930dd55e78ffc1469ed844.zh-cn.xlf"), "", "", "984aa78d-5405-411a-8b9f-4541e5a8b93a.38149897969e365477930dd55e78ffc1469ed844.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A8"), (New-MdUrl ".localization-config"), "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "de-de": same shape as zh-cn
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows.Item(6).Insert()
$wsDe.Range("A6").Value2 = $newMdFile
$wsDe.Range("B6").Value2 = $readyStatus
$wsDe.Range("C6").Value2 = $newDeHandoff
$wsDe.Range("D6").Value2 = $deHandoffDate
$wsDe.Range("G6").Value2 = $epoch
$wsDe.Range("H6").Value2 = $includeStatus

$wsDe.Cells.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), (New-MdUrl "20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.md"), "", "", "20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), (New-DeHandoffUrl "20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.550666d79c2cc553a731e154650f5b8682684363.de-de.xlf"), "", "", "20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.550666d79c2cc553a731e154650f5b8682684363.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), (New-MdUrl "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md"), "", "", "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), (New-DeHandoffUrl "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.de-de.xlf"), "", "", "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), (New-MdUrl "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md"), "", "", "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), (New-DeHandoffUrl "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.de-de.xlf"), "", "", "3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), (New-MdUrl "831315f2-da6a-4fc4-b117-60389f3c6074.md"), "", "", "831315f2-da6a-4fc4-b117-60389f3c6074.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), (New-DeHandoffUrl "831315f2-da6a-4fc4-b117-60389f3c6074.33096a8b89da196732c33acf66abed619fc16792.de-de.xlf"), "", "", "831315f2-da6a-4fc4-b117-60389f3c6074.33096a8b89da196732c33acf66abed619fc16792.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), (New-MdUrl "ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.md"), "", "", "ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C5"), (New-DeHandoffUrl "ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.fc1d393ccc271987bd0fbbef8d59cb1dfbb601f9.de-de.xlf"), "", "", "ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.fc1d393ccc271987bd0fbbef8d59cb1dfbb601f9.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), (New-MdUrl $newMdFile), "", "", $newMdFile)
$wsDe.Hyperlinks.Add($wsDe.Range("C6"), (New-DeHandoffUrl $newDeHandoff), "", "", $newDeHandoff)
$wsDe.Hyperlinks.Add($wsDe.Range("A7"), (New-MdUrl "984aa78d-5405-411a-8b9f-4541e5a8b93a.md"), "", "", "984aa78d-5405-411a-8b9f-4541e5a8b93a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C7"), (New-DeHandoffUrl "984aa78d-5405-411a-8b9f-4541e5a8b93a.38149897969e365477930dd55e78ffc1469ed844.de-de.xlf"), "", "", "984aa78d-5405-411a-8b9f-4541e5a8b93a.38149897969e365477930dd55e78ffc1469ed844.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A8"), (New-MdUrl ".localization-config"), "", "", ".localization-config")

Write-Host "Report generated for handoff of $newMdFile"
